{"js": "// Applies the \"Added many more features\" revision to the Black Widow\n// slot game review: retitles the page/meta title, rewrites the\n// \"What we like\" / \"What we don't like\" bullet lists, and tightens the\n// closing meta-description paragraph.\n\nconst replacements = [\n  // Title (appears twice: the H1 heading and the bolded SEO title run\n  // near the end of the document) \u2014 both instances get the same text.\n  {\n    find: \"Play Black Widow Slot Game for Free | Review\",\n    replace: \"Play Black Widow Slot Game Free\",\n  },\n  // \"What we like\" bullets\n  {\n    find: \"Thrilling game theme based on a popular movie\",\n    replace: \"Thrilling gameplay based on a popular movie\",\n  },\n  {\n    find: \"High payout for landing five Black Widow symbols\",\n    replace: \"Mysterious graphics with a fitting theme\",\n  },\n  {\n    find: \"Free spins bonus round with up to 98 free spins\",\n    replace: \"Chance to win huge prizes\",\n  },\n  {\n    find: \"Available on both mobile and desktop devices\",\n    replace: \"Free spins bonus round\",\n  },\n  // \"What we don't like\" bullets\n  {\n    find: \"Lower payout values for landing male symbols on reel three\",\n    replace: \"Limited variety of bonus rounds\",\n  },\n  {\n    find: \"Limited spider-themed slot options\",\n    replace: \"No progressive jackpot\",\n  },\n  // Closing meta-description paragraph\n  {\n    find:\n      \"Read our review of Black Widow slot game, available to play for free. Get a chance to win huge prizes with the game's free spins bonus round.\",\n    replace:\n      \"Read our review of Black Widow slot game and play it for free. Win huge prizes!\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the \"Added many more features\" revision to the Black Widow\n# slot game review: retitles the page/meta title, rewrites the\n# \"What we like\" / \"What we don't like\" bullet lists, and tightens the\n# closing meta-description paragraph.\n#\n# Each target string is the full (trimmed) text of a single paragraph in\n# this document, so paragraphs are matched by exact text and edited via\n# Find/Replace scoped to that paragraph's own Range -- this keeps the\n# same-looking phrase elsewhere in the review (e.g. the lower-case\n# \"available on both mobile and desktop devices\" used inside two other\n# sentences) untouched.\n#\n# Word's paragraph-rewrite-on-replace silently drops a leading empty\n# <w:r/> run whenever the paragraph also carries a <w:pPr> (several of the\n# bullet paragraphs below have exactly that shape). Restore-LeadingEmptyRun\n# puts that empty run back so paragraph structure is otherwise undisturbed.\n\nfunction Get-ParagraphBodyFragment($doc, $paraIndex) {\n    # Range.WordOpenXML always returns a full pkg:package (styles, theme,\n    # etc. included); slice out just the <w:body>...</w:body> piece that\n    # describes this paragraph so the checks below aren't confused by\n    # unrelated <w:pPr>/<w:r> elements elsewhere in the package.\n    $p = $doc.Paragraphs.Item($paraIndex)\n    $xml = $p.Range.WordOpenXML\n    $bodyStart = $xml.IndexOf(\"<w:body>\")\n    $bodyEnd = $xml.IndexOf(\"</w:body>\")\n    return $xml.Substring($bodyStart, $bodyEnd - $bodyStart)\n}\n\nfunction Test-HasLeadingEmptyRun($doc, $paraIndex) {\n    $frag = Get-ParagraphBodyFragment $doc $paraIndex\n    return ($frag -match '<w:r\\s*/>\\s*<w:r[\\s>]' -or $frag -match '<w:r>\\s*</w:r>\\s*<w:r[\\s>]')\n}\n\nfunction Restore-LeadingEmptyRun($doc, $paraIndex) {\n    $p = $doc.Paragraphs.Item($paraIndex)\n    $insertPoint = $doc.Range($p.Range.Start, $p.Range.Start)\n    $frag = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $insertPoint.InsertXML($frag)\n}\n\nfunction Find-ParagraphIndexesByText($doc, $text) {\n    $paras = $doc.Paragraphs\n    $count = $paras.Count\n    $result = @()\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $paras.Item($i)\n        $t = $p.Range.Text.TrimEnd([char]13)\n        if ($t -eq $text) {\n            $result += $i\n        }\n    }\n    return $result\n}\n\nfunction Replace-ParagraphText($doc, $findText, $replaceText) {\n    $indexes = Find-ParagraphIndexesByText $doc $findText\n    foreach ($idx in $indexes) {\n        # Figure out (before editing) whether this paragraph is at risk of\n        # losing its leading empty run during the text replace.\n        $hasPPr = (Get-ParagraphBodyFragment $doc $idx) -match '<w:pPr>'\n        $hadEmptyRun = Test-HasLeadingEmptyRun $doc $idx\n        $needsRestore = $hasPPr -and $hadEmptyRun\n\n        $rng = $doc.Paragraphs.Item($idx).Range\n        $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n\n        if ($needsRestore) {\n            Restore-LeadingEmptyRun $doc $idx\n        }\n    }\n}\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    # Title (appears twice: the H1 heading and the bolded SEO title run\n    # near the end of the document) \u2014 both instances get the same text.\n    @(\"Play Black Widow Slot Game for Free | Review\", \"Play Black Widow Slot Game Free\"),\n\n    # \"What we like\" bullets\n    @(\"Thrilling game theme based on a popular movie\", \"Thrilling gameplay based on a popular movie\"),\n    @(\"High payout for landing five Black Widow symbols\", \"Mysterious graphics with a fitting theme\"),\n    @(\"Free spins bonus round with up to 98 free spins\", \"Chance to win huge prizes\"),\n    @(\"Available on both mobile and desktop devices\", \"Free spins bonus round\"),\n\n    # \"What we don't like\" bullets\n    @(\"Lower payout values for landing male symbols on reel three\", \"Limited variety of bonus rounds\"),\n    @(\"Limited spider-themed slot options\", \"No progressive jackpot\"),\n\n    # Closing meta-description paragraph\n    @(\"Read our review of Black Widow slot game, available to play for free. Get a chance to win huge prizes with the game's free spins bonus round.\", \"Read our review of Black Widow slot game and play it for free. Win huge prizes!\")\n)\n\nforeach ($pair in $replacements) {\n    Replace-ParagraphText $d $pair[0] $pair[1]\n}\n"}
